$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1202.185677818615
$ws.Range("C2").Value = -1202.185677818615
$ws.Range("D2").Value = -1202.185677818615

$ws.Range("B3").Value = -8.438828532583759
$ws.Range("C3").Value = -6.792568641035932
$ws.Range("D3").Value = -107.8213689204281

$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.790683654242312
$ws.Range("C4").Value = 0.8043968418901722
$ws.Range("D4").Value = -24.36813127761268

$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = -3.988025472442858
$ws.Range("C5").Value = -6.875104746379869
$ws.Range("D5").Value = -518.1470808205112
